$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the placeholder chart ("Chart 2") that plotted the raw B16:B19 /
# C16:C19 helper values -- that data + chart got moved/cleaned up as part
# of relocating this file into its folder. The other chart ("Chart 3",
# Silicone Rubber / PVC series) is left untouched.
for ($i = $ws.ChartObjects().Count; $i -ge 1; $i--) {
    $co = $ws.ChartObjects($i)
    if ($co.Name -eq "Chart 2") {
        $co.Delete()
    }
}

# Clear the now-unused helper values (keep their formatting/style).
$ws.Range("B16:C19").ClearContents()

# Restore the view state (scroll position + active selection) recorded
# after the edit.
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G20").Select()
